# Fruta / hortaliza, semanal
#
# A new weekly price-report row is inserted at row 336 (pushing the
# existing rows 336-355 down to 337-356). The new row carries the same
# market / product / variety / grade / unit metadata as the row that used
# to sit at 336, but with a new sampling date and new price figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 336:355 down to 337:356, leaving a blank row 336 behind with
# the formatting of the row above copied down (matches Excel's own
# "Insert" behaviour for a whole row).
$ws.Rows.Item(336).EntireRow.Insert()

# Populate the newly inserted row 336.
$ws.Cells.Item(336, 1).Value = 4
$ws.Cells.Item(336, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(336, 3).Value = "Los Lagos"
$ws.Cells.Item(336, 4).Value2 = 44931
$ws.Cells.Item(336, 5).Value = 10
$ws.Cells.Item(336, 6).Value = 100112021
$ws.Cells.Item(336, 7).Value = "Ají"
$ws.Cells.Item(336, 8).Value = "Inferno"
$ws.Cells.Item(336, 9).Value = "Primera"
$ws.Cells.Item(336, 10).Value = 80
$ws.Cells.Item(336, 11).Value = 24000
$ws.Cells.Item(336, 12).Value = 24000
$ws.Cells.Item(336, 13).Value = 24000
$ws.Cells.Item(336, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(336, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(336, 16).Value = 2400
$ws.Cells.Item(336, 17).Value = 10
$ws.Cells.Item(336, 18).Value = "Hortaliza"
